$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp (A) and load (B) values for rows 2-93
$ws.Range("A2").Value = 45526
$ws.Range("B2").Value = 5580
$ws.Range("A3").Value = 45526.01041666666
$ws.Range("B3").Value = 5510
$ws.Range("A4").Value = 45526.02083333334
$ws.Range("B4").Value = 5450
$ws.Range("A5").Value = 45526.03125
$ws.Range("B5").Value = 5400
$ws.Range("A6").Value = 45526.04166666666
$ws.Range("B6").Value = 5360
$ws.Range("A7").Value = 45526.05208333334
$ws.Range("B7").Value = 5320
$ws.Range("A8").Value = 45526.0625
$ws.Range("B8").Value = 5290
$ws.Range("A9").Value = 45526.07291666666
$ws.Range("B9").Value = 5260
$ws.Range("A10").Value = 45526.08333333334
$ws.Range("B10").Value = 5240
$ws.Range("A11").Value = 45526.09375
$ws.Range("B11").Value = 5220
$ws.Range("A12").Value = 45526.10416666666
$ws.Range("B12").Value = 5210
$ws.Range("A13").Value = 45526.11458333334
$ws.Range("B13").Value = 5200
$ws.Range("A14").Value = 45526.125
$ws.Range("B14").Value = 5200
$ws.Range("A15").Value = 45526.13541666666
$ws.Range("B15").Value = 5200
$ws.Range("A16").Value = 45526.14583333334
$ws.Range("B16").Value = 5210
$ws.Range("A17").Value = 45526.15625
$ws.Range("B17").Value = 5220
$ws.Range("A18").Value = 45526.16666666666
$ws.Range("B18").Value = 5240
$ws.Range("A19").Value = 45526.17708333334
$ws.Range("B19").Value = 5270
$ws.Range("A20").Value = 45526.1875
$ws.Range("B20").Value = 5300
$ws.Range("A21").Value = 45526.19791666666
$ws.Range("B21").Value = 5360
$ws.Range("A22").Value = 45526.20833333334
$ws.Range("B22").Value = 5420
$ws.Range("A23").Value = 45526.21875
$ws.Range("B23").Value = 5500
$ws.Range("A24").Value = 45526.22916666666
$ws.Range("B24").Value = 5600
$ws.Range("A25").Value = 45526.23958333334
$ws.Range("B25").Value = 5710
$ws.Range("A26").Value = 45526.25
$ws.Range("B26").Value = 5830
$ws.Range("A27").Value = 45526.26041666666
$ws.Range("B27").Value = 5960
$ws.Range("A28").Value = 45526.27083333334
$ws.Range("B28").Value = 6080
$ws.Range("A29").Value = 45526.28125
$ws.Range("B29").Value = 6190
$ws.Range("A30").Value = 45526.29166666666
$ws.Range("B30").Value = 6280
$ws.Range("A31").Value = 45526.30208333334
$ws.Range("B31").Value = 6350
$ws.Range("A32").Value = 45526.3125
$ws.Range("B32").Value = 6400
$ws.Range("A33").Value = 45526.32291666666
$ws.Range("B33").Value = 6400
$ws.Range("A34").Value = 45526.33333333334
$ws.Range("B34").Value = 6400
$ws.Range("A35").Value = 45526.34375
$ws.Range("B35").Value = 6390
$ws.Range("A36").Value = 45526.35416666666
$ws.Range("B36").Value = 6350
$ws.Range("A37").Value = 45526.36458333334
$ws.Range("B37").Value = 6310
$ws.Range("A38").Value = 45526.375
$ws.Range("B38").Value = 6260
$ws.Range("A39").Value = 45526.38541666666
$ws.Range("B39").Value = 6220
$ws.Range("A40").Value = 45526.39583333334
$ws.Range("B40").Value = 6190
$ws.Range("A41").Value = 45526.40625
$ws.Range("B41").Value = 6170
$ws.Range("A42").Value = 45526.41666666666
$ws.Range("B42").Value = 6160
$ws.Range("A43").Value = 45526.42708333334
$ws.Range("B43").Value = 6160
$ws.Range("A44").Value = 45526.4375
$ws.Range("B44").Value = 6170
$ws.Range("A45").Value = 45526.44791666666
$ws.Range("B45").Value = 6180
$ws.Range("A46").Value = 45526.45833333334
$ws.Range("B46").Value = 6180
$ws.Range("A47").Value = 45526.46875
$ws.Range("B47").Value = 6190
$ws.Range("A48").Value = 45526.47916666666
$ws.Range("B48").Value = 6190
$ws.Range("A49").Value = 45526.48958333334
$ws.Range("B49").Value = 6190
$ws.Range("A50").Value = 45526.5
$ws.Range("B50").Value = 6190
$ws.Range("A51").Value = 45526.51041666666
$ws.Range("B51").Value = 6200
$ws.Range("A52").Value = 45526.52083333334
$ws.Range("B52").Value = 6220
$ws.Range("A53").Value = 45526.53125
$ws.Range("B53").Value = 6240
$ws.Range("A54").Value = 45526.54166666666
$ws.Range("B54").Value = 6270
$ws.Range("A55").Value = 45526.55208333334
$ws.Range("B55").Value = 6300
$ws.Range("A56").Value = 45526.5625
$ws.Range("B56").Value = 6330
$ws.Range("A57").Value = 45526.57291666666
$ws.Range("B57").Value = 6360
$ws.Range("A58").Value = 45526.58333333334
$ws.Range("B58").Value = 6380
$ws.Range("A59").Value = 45526.59375
$ws.Range("B59").Value = 6400
$ws.Range("A60").Value = 45526.60416666666
$ws.Range("B60").Value = 6430
$ws.Range("A61").Value = 45526.61458333334
$ws.Range("B61").Value = 6460
$ws.Range("A62").Value = 45526.625
$ws.Range("B62").Value = 6500
$ws.Range("A63").Value = 45526.63541666666
$ws.Range("B63").Value = 6560
$ws.Range("A64").Value = 45526.64583333334
$ws.Range("B64").Value = 6630
$ws.Range("A65").Value = 45526.65625
$ws.Range("B65").Value = 6710
$ws.Range("A66").Value = 45526.66666666666
$ws.Range("B66").Value = 6800
$ws.Range("A67").Value = 45526.67708333334
$ws.Range("B67").Value = 6890
$ws.Range("A68").Value = 45526.6875
$ws.Range("B68").Value = 6970
$ws.Range("A69").Value = 45526.69791666666
$ws.Range("B69").Value = 7030
$ws.Range("A70").Value = 45526.70833333334
$ws.Range("B70").Value = 7090
$ws.Range("A71").Value = 45526.71875
$ws.Range("B71").Value = 7140
$ws.Range("A72").Value = 45526.72916666666
$ws.Range("B72").Value = 7200
$ws.Range("A73").Value = 45526.73958333334
$ws.Range("B73").Value = 7270
$ws.Range("A74").Value = 45526.75
$ws.Range("B74").Value = 7350
$ws.Range("A75").Value = 45526.76041666666
$ws.Range("B75").Value = 7410
$ws.Range("A76").Value = 45526.77083333334
$ws.Range("B76").Value = 7460
$ws.Range("A77").Value = 45526.78125
$ws.Range("B77").Value = 7470
$ws.Range("A78").Value = 45526.79166666666
$ws.Range("B78").Value = 7480
$ws.Range("A79").Value = 45526.80208333334
$ws.Range("B79").Value = 7500
$ws.Range("A80").Value = 45526.8125
$ws.Range("B80").Value = 7500
$ws.Range("A81").Value = 45526.82291666666
$ws.Range("B81").Value = 7500
$ws.Range("A82").Value = 45526.83333333334
$ws.Range("B82").Value = 7500
$ws.Range("A83").Value = 45526.84375
$ws.Range("B83").Value = 7500
$ws.Range("A84").Value = 45526.85416666666
$ws.Range("B84").Value = 7460
$ws.Range("A85").Value = 45526.86458333334
$ws.Range("B85").Value = 7320
$ws.Range("A86").Value = 45526.875
$ws.Range("B86").Value = 7120
$ws.Range("A87").Value = 45526.88541666666
$ws.Range("B87").Value = 6980
$ws.Range("A88").Value = 45526.89583333334
$ws.Range("B88").Value = 6820
$ws.Range("A89").Value = 45526.90625
$ws.Range("B89").Value = 6670
$ws.Range("A90").Value = 45526.91666666666
$ws.Range("B90").Value = 6440
$ws.Range("A91").Value = 45526.92708333334
$ws.Range("B91").Value = 6260
$ws.Range("A92").Value = 45526.9375
$ws.Range("B92").Value = 6180
$ws.Range("A93").Value = 45526.94791666666
$ws.Range("B93").Value = 6070

# Remove now-unused rows 94-98 (shrinks dimension to A1:B93)
$ws.Range("A94:A98").EntireRow.Delete()
